# Generate Report for Handoff
# Updates the "Latest Handoff Date/Datetime" timestamps for the last
# (be32e1f2-...) file row on the Overview, zh-cn and de-de sheets, as
# part of regenerating the localization-status handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date", row 7 = be32e1f2-...
$overview.Range("D7").Value = "2016-03-23 14:42:12"

# zh-cn sheet: column E = "Latest Handoff Datetime", row 7 = be32e1f2-...
$zhcn.Range("E7").Value = "2016-03-23 14:42:08"

# de-de sheet: column E = "Latest Handoff Datetime", row 7 = be32e1f2-...
$dede.Range("E7").Value = "2016-03-23 14:42:12"
